$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (34-36) of master data, following the existing
# pattern in the sheet (regcntr_id, usr_id, machine_id, lang_code,
# is_active, cr_by, cr_dtimes).
$newRows = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the trailing selection state left behind in the saved workbook:
# entire rows selected starting just after the new data (row 37) through
# the bottom of the sheet.
$ws.Rows("37:1048576").Select()
